{"js": "// Update Xbox GDK Samples to November GDK release:\n// \"If using Project Scarlett, set the active solution platform to \" ->\n// \"If using an Xbox One X|S devkit, set the active solution platform to \"\n\nconst searchResults = context.document.body.search(\"If using Project Scarlett, set the active solution platform to \", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\n    \"If using an Xbox One X|S devkit, set the active solution platform to \",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Update Xbox GDK Samples to November GDK release:\n# \"If using Project Scarlett, set the active solution platform to \" ->\n# \"If using an Xbox One X|S devkit, set the active solution platform to \"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"If using Project Scarlett, set the active solution platform to \"\n$find.Replacement.Text = \"If using an Xbox One X|S devkit, set the active solution platform to \"\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n#              Format, ReplaceWith, Replace:=wdReplaceOne)\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n"}
